$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 73
$ws1.Range("F6").Value = 122
$ws1.Range("F7").Value = 1202
$ws1.Range("F8").Value = 1495
$ws1.Range("F10").Value = 372
$ws1.Range("F12").Value = 131
$ws1.Range("F16").Value = 267
$ws1.Range("F19").Value = 1699
$ws1.Range("F21").Value = 105
$ws1.Range("F23").Value = 645
$ws1.Range("F25").Value = 331
$ws1.Range("F26").Value = 4080
$ws1.Range("F29").Value = 252
$ws1.Range("F31").Value = 128
$ws1.Range("F33").Value = 406
$ws1.Range("F35").Value = 182

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 73
$ws4.Range("F6").Value = 122
$ws4.Range("F7").Value = 1202
$ws4.Range("F8").Value = 1495
$ws4.Range("F10").Value = 372
$ws4.Range("F12").Value = 131
$ws4.Range("F16").Value = 267
$ws4.Range("F19").Value = 1699
$ws4.Range("F21").Value = 105
$ws4.Range("F23").Value = 645
$ws4.Range("F25").Value = 331
$ws4.Range("F26").Value = 4080
$ws4.Range("F29").Value = 252
$ws4.Range("F31").Value = 128
$ws4.Range("F33").Value = 406
$ws4.Range("F35").Value = 182
